{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the last paragraph that holds \"Eigene Texture Layer erstellen\" \u2014\n// the new bullet \"Skybox erstellt\" is added right after it.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"Eigene Texture Layer erstellen\") !== -1) {\n    target = p;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate paragraph 'Eigene Texture Layer erstellen'.\");\n}\n\n// Inserting directly after the existing bullet inherits its paragraph\n// (list/style) and run formatting automatically, matching the new bullet\n// in the diff.\ntarget.insertParagraph(\"Skybox erstellt\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Find the paragraph that holds \"Eigene Texture Layer erstellen\" \u2014 the new\n# bullet \"Skybox erstellt\" belongs right after it in the same bulleted list.\n$targetIndex = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    if ($p.Range.Text -like \"*Eigene Texture Layer erstellen*\") {\n        $targetIndex = $i\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate paragraph 'Eigene Texture Layer erstellen'.\"\n}\n\n$target = $d.Paragraphs.Item($targetIndex)\n\n# Splitting the paragraph mark like this inherits the list/style and run\n# formatting of the source bullet automatically.\n$target.Range.InsertParagraphAfter()\n\n# Re-fetch the freshly created paragraph by its (now stable) index and set\n# its text \u2014 writing through the stale $target.Next() reference doesn't\n# stick, so address the collection directly.\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newPara.Range.Text = \"Skybox erstellt\"\n"}
